$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item("직사각형 26").Delete()
$s.Shapes.Item("직사각형 30").Delete()
$s.Shapes.Item("그래픽 31").Delete()
